$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table (rows 1-6) previously ended at column Q (year 2020). Add a new
# column R for year 2021 by inserting a copy of column Q immediately to its
# right, which brings along the exact formatting used throughout the Q
# column (including the thin "spacer" row 2 and the header/data rows).
$ws.Columns("Q").Copy()
[void]$ws.Columns("R").Insert(-4161, 0)

# Now replace the copied values in row 3 (year), row 4 (GVA share %) and
# row 5 (GVA per capita) with the new 2021 figures. Row 2's single
# formatting-only cell (R2) is left exactly as copied (blank, same style
# as Q2).
$ws.Range("R3").Value = 2021
$ws.Range("R4").Value = 13.5
$ws.Range("R5").Value = 15.1

# Move the active selection to reflect the new state of the sheet.
[void]$ws.Range("T3").Select()
